$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3230.318
$ws.Range("J17").Value = 3868.8286
$ws.Range("L17").Value = 11606.4858
$ws.Range("N17").Value = -11942.4858
$ws.Range("H70").Value = 70533.13
$ws.Range("J70").Value = 4545.5454
$ws.Range("L70").Value = 13636.6362
$ws.Range("N70").Value = -14176.6362
$ws.Range("H73").Value = 70533.13
$ws.Range("J73").Value = 4545.5454
$ws.Range("L73").Value = 13636.6362
$ws.Range("N73").Value = -15508.6362
$ws.Range("H98").Value = 4430.069
$ws.Range("I98").Value = 999
$ws.Range("J98").Value = 13436.625
$ws.Range("K98").Value = 999
$ws.Range("L98").Value = 13436.625
$ws.Range("M98").Value = 499
$ws.Range("N98").Value = -16432.625
$ws.Range("H122").Value = 4430.069
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 13436.625
$ws.Range("K122").Value = 2997
$ws.Range("L122").Value = 40309.875
$ws.Range("M122").Value = -547
$ws.Range("N122").Value = -45209.875
$ws.Range("H133").Value = 60408.43
$ws.Range("J133").Value = 60408.43
$ws.Range("L133").Value = 60408.43
$ws.Range("N133").Value = -70528.42999999999
$ws.Range("H138").Value = 5973
$ws.Range("I138").Value = 3223.3333
$ws.Range("J138").Value = 6972.879
$ws.Range("K138").Value = 9669.999899999999
$ws.Range("L138").Value = 20918.637
$ws.Range("M138").Value = -4529.999899999999
$ws.Range("N138").Value = -31198.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26427.75
$ws.Range("I2").Value = 355.5
$ws.Range("J2").Value = 52500
$ws.Range("K2").Value = 355.5
$ws.Range("L2").Value = 52500
$ws.Range("M2").Value = -242.5
$ws.Range("N2").Value = -52726
$ws.Range("H88").Value = 4195
$ws.Range("J88").Value = 3310.5
$ws.Range("L88").Value = 3310.5
$ws.Range("N88").Value = -4122.5
$ws.Range("H91").Value = 4195
$ws.Range("J91").Value = 3310.5
$ws.Range("L91").Value = 3310.5
$ws.Range("N91").Value = -6118.5
$ws.Range("H116").Value = 26427.75
$ws.Range("I116").Value = 355.5
$ws.Range("J116").Value = 52500
$ws.Range("K116").Value = 355.5
$ws.Range("L116").Value = 52500
$ws.Range("M116").Value = 1938.5
$ws.Range("N116").Value = -57088
$ws.Range("H119").Value = 65152.8
$ws.Range("J119").Value = 65152.8
$ws.Range("L119").Value = 65152.8
$ws.Range("N119").Value = -74828.8
$ws.Range("H139").Value = 59849.25
$ws.Range("J139").Value = 59849.25
$ws.Range("L139").Value = 59849.25
$ws.Range("N139").Value = -70129.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26427.75
$ws.Range("I3").Value = 355.5
$ws.Range("J3").Value = 52500
$ws.Range("K3").Value = 355.5
$ws.Range("L3").Value = 52500
$ws.Range("M3").Value = -241.5
$ws.Range("N3").Value = -52728

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 52342.906
$ws.Range("J31").Value = 78227.36
$ws.Range("L31").Value = 78227.36
$ws.Range("N31").Value = -78817.36
$ws.Range("H34").Value = 52342.906
$ws.Range("J34").Value = 78227.36
$ws.Range("L34").Value = 78227.36
$ws.Range("N34").Value = -78631.36
$ws.Range("H58").Value = 8254.375
$ws.Range("I58").Value = 2103.6667
$ws.Range("J58").Value = 11944.8
$ws.Range("K58").Value = 2103.6667
$ws.Range("L58").Value = 11944.8
$ws.Range("M58").Value = -1900.6667
$ws.Range("N58").Value = -12350.8
$ws.Range("H62").Value = 10190.462
$ws.Range("I62").Value = 3632.8
$ws.Range("J62").Value = 14289
$ws.Range("K62").Value = 3632.8
$ws.Range("L62").Value = 14289
$ws.Range("M62").Value = -3008.8
$ws.Range("N62").Value = -15537
$ws.Range("H65").Value = 10190.462
$ws.Range("I65").Value = 3632.8
$ws.Range("J65").Value = 14289
$ws.Range("K65").Value = 18164
$ws.Range("L65").Value = 71445
$ws.Range("M65").Value = -15044
$ws.Range("N65").Value = -77685
$ws.Range("H125").Value = 99675
$ws.Range("J125").Value = 99675
$ws.Range("L125").Value = 99675
$ws.Range("N125").Value = -104595
$ws.Range("H136").Value = 8254.375
$ws.Range("I136").Value = 2103.6667
$ws.Range("J136").Value = 11944.8
$ws.Range("K136").Value = 6311.000100000001
$ws.Range("L136").Value = 35834.39999999999
$ws.Range("M136").Value = -3761.000100000001
$ws.Range("N136").Value = -40934.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 21.61111
$ws.Range("J12").Value = 11.6875
$ws.Range("L12").Value = 35.0625
$ws.Range("N12").Value = -381.0625
$ws.Range("H64").Value = 90915440
$ws.Range("I64").Value = 142861060
$ws.Range("K64").Value = 428583180
$ws.Range("M64").Value = -428582910
$ws.Range("H67").Value = 90915440
$ws.Range("I67").Value = 142861060
$ws.Range("K67").Value = 428583180
$ws.Range("M67").Value = -428582244
$ws.Range("H120").Value = 20000
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("H139").Value = 8242.6
$ws.Range("I139").Value = 3000
$ws.Range("J139").Value = 11737.667
$ws.Range("K139").Value = 9000
$ws.Range("L139").Value = 35213.001
$ws.Range("M139").Value = -3860
$ws.Range("N139").Value = -45493.001
$ws.Range("M120").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 553.8889
$ws.Range("J2").Value = 1065.8889
$ws.Range("L2").Value = 1065.8889
$ws.Range("N2").Value = -1291.8889
$ws.Range("H10").Value = 15668.667
$ws.Range("I10").Value = 15668.667
$ws.Range("K10").Value = 15668.667
$ws.Range("M10").Value = -15499.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7565.077
$ws.Range("I7").Value = 4420.8423
$ws.Range("J7").Value = 16099.429
$ws.Range("K7").Value = 4420.8423
$ws.Range("L7").Value = 16099.429
$ws.Range("M7").Value = -4308.8423
$ws.Range("N7").Value = -16323.429
$ws.Range("H93").Value = 100000
$ws.Range("I93").Value = 100000
$ws.Range("K93").Value = 100000
$ws.Range("H99").Value = 29000
$ws.Range("I99").Value = 29000
$ws.Range("K99").Value = 29000
$ws.Range("H126").Value = 7565.077
$ws.Range("I126").Value = 4420.8423
$ws.Range("J126").Value = 16099.429
$ws.Range("K126").Value = 13262.5269
$ws.Range("L126").Value = 48298.287
$ws.Range("M126").Value = -10792.5269
$ws.Range("N126").Value = -53238.287
$ws.Range("M93").Value = -98752
$ws.Range("M99").Value = -26005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4059.375
$ws.Range("I136").Value = 2115.238
$ws.Range("J136").Value = 17668.334
$ws.Range("K136").Value = 6345.714
$ws.Range("L136").Value = 53005.00199999999
$ws.Range("M136").Value = -3795.714
$ws.Range("N136").Value = -58105.00199999999
